$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The biosampleNumber column (C) is renumbered: the new CBF1 expb samples
# continue the existing sequence, shifting 1..12 up to 21..32.
for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value2 = $cell.Value2 + 20
}

# The renumbered biosampleNumber column now gets a distinct look: Calibri
# 11pt black, right-aligned (it previously inherited the plain Arial style).
$rng = $ws.Range("C2:C13")
$rng.Font.Color = 0
$rng.Font.Size = 11
$rng.Font.Name = "Calibri"
$rng.HorizontalAlignment = -4152

# Restore the user's last selection / active cell.
[void]$ws.Range("D24").Select()
